$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D2:D51 and E2:E51 are treated as Text so numeric-looking
# strings (e.g. "0.9994", "26.29") are not auto-converted to numbers.
$ws.Range("D2:D51,E2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "29.945.64"
$ws.Range("E2").Value = "  +0.46%  "
$ws.Range("D3").Value = "1.907.66"
$ws.Range("E3").Value = "  +0.80%  "
$ws.Range("D4").Value = "0.9994"
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").Value = "0.8006"
$ws.Range("E5").Value = "  +5.83%  "
$ws.Range("D6").Value = "241.60"
$ws.Range("E6").Value = "  +1.08%  "
$ws.Range("D7").Value = "0.9991"
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("D8").Value = "0.3155"
$ws.Range("E8").Value = "  +3.79%  "
$ws.Range("D9").Value = "26.29"
$ws.Range("E9").Value = "  +3.96%  "
$ws.Range("D10").Value = "0.06906"
$ws.Range("E10").Value = "  +1.25%  "
$ws.Range("D11").Value = "0.07986"
$ws.Range("E11").Value = "  +0.10%  "
$ws.Range("D12").Value = "1.909.29"
$ws.Range("E12").Value = "  +0.91%  "
$ws.Range("D13").Value = "0.7369"
$ws.Range("E13").Value = "  -1.46%  "
$ws.Range("D14").Value = "5.187"
$ws.Range("E14").Value = "  -0.14%  "
$ws.Range("D15").Value = "92.87"
$ws.Range("E15").Value = "  +1.99%  "
$ws.Range("D16").Value = "29.936.82"
$ws.Range("E16").Value = "  +0.43%  "
$ws.Range("D17").Value = "13.95"
$ws.Range("E17").Value = "  +0.58%  "
$ws.Range("D18").Value = "5.860"
$ws.Range("E18").Value = "  -2.23%  "
$ws.Range("D19").Value = "245.25"
$ws.Range("E19").Value = "  +4.62%  "
$ws.Range("D20").Value = "0.000007733"
$ws.Range("E20").Value = "  +0.87%  "
$ws.Range("D21").Value = "0.9997"
$ws.Range("E21").Value = "  -0.10%  "
$ws.Range("D22").Value = "2.147.74"
$ws.Range("E22").Value = "  +0.13%  "
$ws.Range("D23").Value = "0.9997"
$ws.Range("E23").Value = "  -0.17%  "
$ws.Range("D24").Value = "6.816"
$ws.Range("E24").Value = "  -1.47%  "
$ws.Range("D25").Value = "167.61"
$ws.Range("E25").Value = "  +1.36%  "
$ws.Range("D26").Value = "9.187"
$ws.Range("E26").Value = "  -0.65%  "
$ws.Range("D27").Value = "0.1411"
$ws.Range("E27").Value = "  +9.77%  "
$ws.Range("D28").Value = "18.89"
$ws.Range("E28").Value = "  +1.20%  "
$ws.Range("D29").Value = "2.028"
$ws.Range("E29").Value = "  -1.42%  "
$ws.Range("D30").Value = "1.361"
$ws.Range("E30").Value = "  +1.60%  "
$ws.Range("D32").Value = "4.299"
$ws.Range("E32").Value = "  +0.52%  "
$ws.Range("D33").Value = "4.081"
$ws.Range("E33").Value = "  +1.49%  "
$ws.Range("D34").Value = "0.05459"
$ws.Range("E34").Value = "  +1.98%  "
$ws.Range("D35").Value = "1.261"
$ws.Range("E35").Value = "  +1.49%  "
$ws.Range("D36").Value = "0.7284"
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("D37").Value = "2.719"
$ws.Range("E37").Value = "  +0.17%  "
$ws.Range("D38").Value = "0.01920"
$ws.Range("E38").Value = "  -0.30%  "
$ws.Range("E39").Value = "  +0.46%  "
$ws.Range("D40").Value = "6.137"
$ws.Range("E40").Value = "  -1.02%  "
$ws.Range("D41").Value = "0.4414"
$ws.Range("E41").Value = "  +0.22%  "
$ws.Range("D42").Value = "72.23"
$ws.Range("E42").Value = "  -0.05%  "
$ws.Range("D43").Value = "0.9995"
$ws.Range("E43").Value = "  -0.10%  "
$ws.Range("D44").Value = "0.8343"
$ws.Range("E44").Value = "  +1.38%  "
$ws.Range("D45").Value = "1.869"
$ws.Range("E45").Value = "  -2.41%  "
$ws.Range("D46").Value = "100.38"
$ws.Range("E46").Value = "  -0.62%  "
$ws.Range("D47").Value = "7.520"
$ws.Range("E47").Value = "  -0.69%  "
$ws.Range("D48").Value = "9.706"
$ws.Range("E48").Value = "  -1.17%  "
$ws.Range("D49").Value = "987.20"
$ws.Range("E49").Value = "  +8.01%  "
$ws.Range("D50").Value = "2.053.76"
$ws.Range("E50").Value = "  +0.40%  "
$ws.Range("D51").Value = "36.16"
$ws.Range("E51").Value = "  +0.39%  "

# Restore default (General) formatting/style so cells match the
# original workbook's unstyled appearance, while keeping the
# values stored as text.
$ws.Range("D2:D51,E2:E51").ClearFormats()

